# Update gh-pages output data for 江西-漫展信息.xlsx
# Sheet 1 = "展览" (rows offset A), Sheet 4 = "全部类型" (has one extra row
# inserted at row 28, so matching events are one row lower from row 29 on).

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (index 1) ----
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(3, 6).Value = 2277      # F3  2276 -> 2277

$ws1.Cells.Item(6, 6).Value = 5161      # F6  5159 -> 5161
$ws1.Cells.Item(6, 7).Value = 70        # G6  65   -> 70

$ws1.Cells.Item(7, 7).Value = 25        # G7  20   -> 25

$ws1.Cells.Item(8, 7).Value = "不可售"   # G8  45 (n) -> "不可售" (str)

$ws1.Cells.Item(9, 6).Value = 316       # F9  314 -> 316

$ws1.Cells.Item(10, 7).Value = 65       # G10 55 -> 65

$ws1.Cells.Item(11, 6).Value = 42       # F11 41 -> 42

$ws1.Cells.Item(12, 6).Value = 221      # F12 220 -> 221

$ws1.Cells.Item(15, 6).Value = 130      # F15 129 -> 130

$ws1.Cells.Item(16, 6).Value = 4229     # F16 4227 -> 4229

$ws1.Cells.Item(17, 6).Value = 762      # F17 761 -> 762

$ws1.Cells.Item(18, 6).Value = 771      # F18 769 -> 771

$ws1.Cells.Item(20, 6).Value = 27       # F20 26 -> 27

$ws1.Cells.Item(29, 6).Value = 1078     # F29 1071 -> 1078

$ws1.Cells.Item(30, 6).Value = 10       # F30 8 -> 10

$ws1.Cells.Item(31, 6).Value = 2707     # F31 2698 -> 2707

$ws1.Cells.Item(32, 6).Value = 434      # F32 432 -> 434

$ws1.Cells.Item(33, 6).Value = 86       # F33 82 -> 86

# ---- Sheet "全部类型" (index 4) ----
# Same underlying events, but this sheet has an extra row (a concert entry)
# at row 28, so every matching row from 29 onward is shifted down by one
# compared to the "展览" sheet.
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(3, 6).Value = 2277      # F3  2276 -> 2277

$ws4.Cells.Item(6, 6).Value = 5161      # F6  5159 -> 5161
$ws4.Cells.Item(6, 7).Value = 70        # G6  65   -> 70

$ws4.Cells.Item(7, 7).Value = 25        # G7  20   -> 25

$ws4.Cells.Item(8, 7).Value = "不可售"   # G8  45 (n) -> "不可售" (str)

$ws4.Cells.Item(9, 6).Value = 316       # F9  314 -> 316

$ws4.Cells.Item(10, 7).Value = 65       # G10 55 -> 65

$ws4.Cells.Item(11, 6).Value = 42       # F11 41 -> 42

$ws4.Cells.Item(12, 6).Value = 221      # F12 220 -> 221

$ws4.Cells.Item(15, 6).Value = 130      # F15 129 -> 130

$ws4.Cells.Item(16, 6).Value = 4229     # F16 4227 -> 4229

$ws4.Cells.Item(17, 6).Value = 762      # F17 761 -> 762

$ws4.Cells.Item(18, 6).Value = 771      # F18 769 -> 771

$ws4.Cells.Item(20, 6).Value = 27       # F20 26 -> 27

$ws4.Cells.Item(30, 6).Value = 1078     # F30 1071 -> 1078

$ws4.Cells.Item(31, 6).Value = 10       # F31 8 -> 10

$ws4.Cells.Item(32, 6).Value = 2707     # F32 2698 -> 2707

$ws4.Cells.Item(33, 6).Value = 434      # F33 432 -> 434

$ws4.Cells.Item(34, 6).Value = 86       # F34 82 -> 86
